$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Helper: set the text of paragraph $idx (1-based) without touching the
# paragraph mark, then (optionally) reapply run formatting on the new
# text range. Returns nothing; operates in place on $d.
# ---------------------------------------------------------------------
function Set-ParaText($idx, $newText) {
    $p = $d.Paragraphs.Item($idx)
    $r = $p.Range
    $r2 = $d.Range($r.Start, $r.End - 1)
    $r2.Text = $newText
}

function Set-ParaFormat($idx, $alignment, $bold, $color, $size) {
    $p = $d.Paragraphs.Item($idx)
    if ($alignment -ne $null) {
        $p.Alignment = $alignment
    }
    $r = $p.Range
    $r2 = $d.Range($r.Start, $r.End - 1)
    if ($bold -ne $null) {
        $r2.Font.Bold = $bold
    }
    if ($color -ne $null) {
        $r2.Font.Color = $color
    }
    if ($size -ne $null) {
        $r2.Font.Size = $size
    }
}

$BLUE = 16711680

# -----------------------------------------------------------------
# 1. Title: "The Great Renovation" -> "Basement Renovation Project",
#    size 48 -> 36 (half-points 96 -> 72, i.e. pt 24 -> 18... actually
#    sz is half-points directly: 48 half-pts = 24pt -> 36 half-pts = 18pt)
# -----------------------------------------------------------------
Set-ParaText 1 "Basement Renovation Project"
Set-ParaFormat 1 $null $null $null 18

# -----------------------------------------------------------------
# 2. Insert new paragraph after paragraph 1: "Location: ..."
# -----------------------------------------------------------------
$p1 = $d.Paragraphs.Item(1)
$p1.Range.InsertParagraphAfter()
$newLoc = $d.Paragraphs.Item(2)
$newLoc.Alignment = 0
$r = $newLoc.Range
$r2 = $d.Range($r.Start, $r.End - 1)
$r2.Text = "Location: 13599 Cobra Dr, Herndon, VA 20171"
$r2.Font.Bold = $false
$r2.Font.Italic = $false
$r2.Font.Color = 0
$r2.Font.Size = 14

# -----------------------------------------------------------------
# 3. "Project Overview" (now paragraph 3): left -> center,
#    color 000000 -> 0000FF, size 36 -> 32 (18pt -> 16pt)
# -----------------------------------------------------------------
Set-ParaFormat 3 1 $null $BLUE 16

# -----------------------------------------------------------------
# 4. "Project Address: ..." (paragraph 4) -> new overview sentence
# -----------------------------------------------------------------
Set-ParaText 4 "This project involves a modern design and premium finishes, with a total labor cost of ($11,250) and is expected to be completed within 2 months."

# -----------------------------------------------------------------
# 5. "Total Estimated Cost: $167,200" (paragraph 5) -> "Cost Breakdown:"
#    color 000000 -> 0000FF (size stays 32/16pt)
# -----------------------------------------------------------------
Set-ParaText 5 "Cost Breakdown:"
Set-ParaFormat 5 $null $null $BLUE $null

# -----------------------------------------------------------------
# 6. Delete "Breakdown of Costs:" (paragraph 6)
# -----------------------------------------------------------------
$d.Paragraphs.Item(6).Range.Delete()

# -----------------------------------------------------------------
# 7. "Labor Costs: $150,000" (now paragraph 6) -> "Labor Cost: ($11,250)"
# -----------------------------------------------------------------
Set-ParaText 6 "Labor Cost: ($11,250)"

# -----------------------------------------------------------------
# 8. "Material Costs: $17,200" (paragraph 7) -> "Material Cost:"
# -----------------------------------------------------------------
Set-ParaText 7 "Material Cost:"

# -----------------------------------------------------------------
# 9. Delete "Material Costs Breakdown:" (paragraph 8)
# -----------------------------------------------------------------
$d.Paragraphs.Item(8).Range.Delete()

# -----------------------------------------------------------------
# 10. "Flooring: ..." (now paragraph 8) -> "  • Flooring: ($4,500)"
# -----------------------------------------------------------------
Set-ParaText 8 "  • Flooring: ($4,500)"

# -----------------------------------------------------------------
# 11. "Baseboard Molding: ..." (paragraph 9) -> "  • Lighting: ($3,000)"
# -----------------------------------------------------------------
Set-ParaText 9 "  • Lighting: ($3,000)"

# -----------------------------------------------------------------
# 12. "Underlayment Material: ..." (paragraph 10) -> "  • Shelving: ($6,500)"
# -----------------------------------------------------------------
Set-ParaText 10 "  • Shelving: ($6,500)"

# -----------------------------------------------------------------
# 13. Delete "Walls and Ceiling Materials:" (paragraph 11)
# -----------------------------------------------------------------
$d.Paragraphs.Item(11).Range.Delete()

# -----------------------------------------------------------------
# 14. "Drywall: ..." (now paragraph 11) -> "  • Electrical and HVAC: ($5,000)"
# -----------------------------------------------------------------
Set-ParaText 11 "  • Electrical and HVAC: ($5,000)"

# -----------------------------------------------------------------
# 15. "Metal Studs: ..." (paragraph 12) -> "  • Walls and Ceiling: ($3,500)"
# -----------------------------------------------------------------
Set-ParaText 12 "  • Walls and Ceiling: ($3,500)"

# -----------------------------------------------------------------
# 16. "Joint Compound: ..." (paragraph 13) -> "  • Games Area: ($7,000)"
# -----------------------------------------------------------------
Set-ParaText 13 "  • Games Area: ($7,000)"

# -----------------------------------------------------------------
# 17. "Paint: ..." (paragraph 14) -> "  • Tinted Glass Wall: ($10,000)"
# -----------------------------------------------------------------
Set-ParaText 14 "  • Tinted Glass Wall: ($10,000)"

# -----------------------------------------------------------------
# 18. Delete "Windows and Doors:" (paragraph 15)
# -----------------------------------------------------------------
$d.Paragraphs.Item(15).Range.Delete()

# -----------------------------------------------------------------
# 19. "Five new, double-hung windows ..." (now paragraph 15) -> "  • Hot Tub: ($15,000)"
# -----------------------------------------------------------------
Set-ParaText 15 "  • Hot Tub: ($15,000)"

# -----------------------------------------------------------------
# 20. "One new, solid-wood entrance door ..." (paragraph 16) -> "  • Miscellaneous: ($14,000)"
# -----------------------------------------------------------------
Set-ParaText 16 "  • Miscellaneous: ($14,000)"

# -----------------------------------------------------------------
# 21. Delete "Electrical Materials:" (paragraph 17)
# -----------------------------------------------------------------
$d.Paragraphs.Item(17).Range.Delete()

# -----------------------------------------------------------------
# 22. "500 feet of 14-gauge wire ..." (now paragraph 17) -> "Total Material Cost: ($68,500)"
# -----------------------------------------------------------------
Set-ParaText 17 "Total Material Cost: ($68,500)"

# -----------------------------------------------------------------
# 23. Delete "20 electrical outlets ..." (paragraph 18) and
#     "10 LED recessed lights ..." (paragraph 18 again, after first delete)
# -----------------------------------------------------------------
$d.Paragraphs.Item(18).Range.Delete()
$d.Paragraphs.Item(18).Range.Delete()

# -----------------------------------------------------------------
# 24. "Plumbing Materials:" (now paragraph 18) -> "Grand Total: ($79,750)"
#     color 000000 -> 0000FF, size 28 -> 32 (14pt -> 16pt)
# -----------------------------------------------------------------
Set-ParaText 18 "Grand Total: ($79,750)"
Set-ParaFormat 18 $null $null $BLUE 16

# -----------------------------------------------------------------
# 25. Delete the next 12 paragraphs (now all at index 19, since each
#     delete shifts everything up):
#       "200 feet of PEX tubing (1/2 inch) ($700)"
#       "15 feet of copper piping (1/2 inch, type M) ($700)"
#       "Two new, low-flow toilets (1.6 GPF) ($1,500)"
#       "One new, single-bowl kitchen sink (stainless steel) ($700)"
#       "Insulation and Drywall Accessories:"
#       "1,500 square feet of fiberglass batt insulation (R-19) ($2,000)"
#       "500 square feet of drywall tape ($700)"
#       "100 drywall screws (coarse-thread) ($700)"
#       "Miscellaneous Items:"
#       "10 pounds of drywall anchors ($200)"
#       "20 pounds of joint compound sandpaper ($200)"
#       "Five rolls of painter's tape ($100)"
# -----------------------------------------------------------------
for ($i = 0; $i -lt 12; $i++) {
    $d.Paragraphs.Item(19).Range.Delete()
}

# -----------------------------------------------------------------
# 26. "Payment Terms:" (now paragraph 19): color 000000 -> 0000FF,
#     size 28 -> 32 (14pt -> 16pt)
# -----------------------------------------------------------------
Set-ParaFormat 19 $null $null $BLUE 16

# -----------------------------------------------------------------
# 27. "Deposit: 20% of total cost ($33,440.00)" (paragraph 20) ->
#     "To initiate the project, a deposit of ($15,950) is required."
# -----------------------------------------------------------------
Set-ParaText 20 "To initiate the project, a deposit of ($15,950) is required."

# -----------------------------------------------------------------
# 28. "Payment Schedule: Monthly payments of $3,344 for 40 months"
#     (paragraph 21) -> "Followed by a payment schedule of $15,950 per
#     month for 4 months, with the final payment due upon completion
#     of the project."
# -----------------------------------------------------------------
Set-ParaText 21 "Followed by a payment schedule of `$15,950 per month for 4 months, with the final payment due upon completion of the project."

# -----------------------------------------------------------------
# 29. Delete "Final Payment: Due upon completion of the project"
#     (paragraph 22)
# -----------------------------------------------------------------
$d.Paragraphs.Item(22).Range.Delete()

Write-Host "Final paragraph count:" $d.Paragraphs.Count
